$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Model Information")
$wsInfo.Range("C3").Value = 75.8

$wsProfile = $wb.Worksheets.Item("Model Profile")
$wsProfile.Range("F4").Value = 20.45356750488281
$wsProfile.Range("G4").Value = 117.1840430034826
$wsProfile.Range("F5").Value = 31.36816978454589
$wsProfile.Range("G5").Value = 156.5747949575072
$wsProfile.Range("F6").Value = 51.50646448135375
$wsProfile.Range("G6").Value = 156.7974428542034
$wsProfile.Range("F7").Value = 100.9795737266541
$wsProfile.Range("G7").Value = 140.3298589991225
$wsProfile.Range("F8").Value = 29.50067996978759
$wsProfile.Range("G8").Value = 78.5649895349803
$wsProfile.Range("F9").Value = 31.62540197372437
$wsProfile.Range("G9").Value = 127.264902264969
$wsProfile.Range("F10").Value = 37.55554676055907
$wsProfile.Range("G10").Value = 176.3530181846559
$wsProfile.Range("F11").Value = 51.32622241973876
$wsProfile.Range("G11").Value = 219.2586291672717
$wsProfile.Range("F12").Value = 19.88508224487304
$wsProfile.Range("G12").Value = 107.7219149462287
$wsProfile.Range("F13").Value = 31.80617332458496
$wsProfile.Range("G13").Value = 141.0505317091839
$wsProfile.Range("F14").Value = 51.76976203918457
$wsProfile.Range("G14").Value = 137.396714248651
$wsProfile.Range("F15").Value = 98.14545869827271
$wsProfile.Range("G15").Value = 133.8823067660375
$wsProfile.Range("F16").Value = 24.84185218811035
$wsProfile.Range("G16").Value = 97.08632510682425
$wsProfile.Range("F17").Value = 28.20687532424926
$wsProfile.Range("G17").Value = 139.4325870695722
$wsProfile.Range("F18").Value = 33.45731496810912
$wsProfile.Range("G18").Value = 171.55825862255
$wsProfile.Range("F19").Value = 50.9054136276245
$wsProfile.Range("G19").Value = 225.5949406713249

$wsRaw = $wb.Worksheets.Item("Model Raw Profile")
$wsRaw.Range("C2").Value = '{"Classifier-resnet50/prepoc-resnet50": {"CPU1": {"THROUGHPUT": [[1, 117.18404300348257], [2, 156.57479495750724], [4, 156.79744285420344], [8, 140.32985899912248]], "LATENCY": [[1, 20.45356750488281], [2, 31.368169784545895], [4, 51.50646448135375], [8, 100.97957372665405]]}}, "Classifier-resnet50/model-resnet50": {"Tesla P40": {"THROUGHPUT": [[1, 78.5649895349803], [2, 127.26490226496902], [4, 176.3530181846559], [8, 219.25862916727175]], "LATENCY": [[1, 29.500679969787594], [2, 31.625401973724365], [4, 37.555546760559075], [8, 51.32622241973876]]}}}'
$wsRaw.Range("C3").Value = '{"Classifier-resnet34/prepoc-resnet34": {"CPU1": {"THROUGHPUT": [[1, 107.72191494622872], [2, 141.0505317091839], [4, 137.39671424865102], [8, 133.88230676603754]], "LATENCY": [[1, 19.885082244873043], [2, 31.806173324584957], [4, 51.76976203918457], [8, 98.1454586982727]]}}, "Classifier-resnet34/model-resnet34": {"Tesla P40": {"THROUGHPUT": [[1, 97.08632510682425], [2, 139.43258706957215], [4, 171.55825862255003], [8, 225.5949406713249]], "LATENCY": [[1, 24.84185218811035], [2, 28.206875324249264], [4, 33.457314968109124], [8, 50.905413627624505]]}}}'
$wsRaw.Range("G3").Value = 75.8
